# Update "想去人数" (interested-people count) figures across sheets to
# reflect freshly generated data (commit: "Update gh-pages to output
# generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# --- 展览 (Exhibitions) sheet ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value  = 12031
$wsExpo.Range("F4").Value  = 36
$wsExpo.Range("F8").Value  = 11918
$wsExpo.Range("F9").Value  = 501
$wsExpo.Range("F10").Value = 1178
$wsExpo.Range("F14").Value = 5902
$wsExpo.Range("F18").Value = 29

# --- 演出 (Performances) sheet ---
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 577

# --- 全部类型 (All types) sheet, aggregates the rows above ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value  = 577
$wsAll.Range("F5").Value  = 12031
$wsAll.Range("F6").Value  = 36
$wsAll.Range("F11").Value = 11918
$wsAll.Range("F12").Value = 501
$wsAll.Range("F13").Value = 1178
$wsAll.Range("F18").Value = 5902
$wsAll.Range("F22").Value = 29
